$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.864.09'
$ws.Range('E2').Value = '  -1.63%  '
$ws.Range('D3').Value = '1.805.56'
$ws.Range('E3').Value = '  -0.98%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '309.90'
$ws.Range('E5').Value = '  -1.44%  '
$ws.Range('E6').Value = '  +0.00%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4658'
$ws.Range('E7').Value = '  +3.73%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3711'
$ws.Range('E8').Value = '  -1.84%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07379'
$ws.Range('E9').Value = '  -0.74%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.8720'
$ws.Range('E10').Value = '  -1.56%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.41'
$ws.Range('E11').Value = '  -2.63%  '
$ws.Range('D12').Value = '1.779.92'
$ws.Range('E12').Value = '  -2.37%  '
$ws.Range('E13').Value = '  -1.41%  '
$ws.Range('E14').Value = '  -0.92%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.477'
$ws.Range('E15').Value = '  -3.77%  '
$ws.Range('E16').Value = '  -1.35%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.002'
$ws.Range('E17').Value = '  +0.01%  '
$ws.Range('E18').Value = '  -1.01%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.68'
$ws.Range('E20').Value = '  -3.04%  '
$ws.Range('D21').Value = '26.859.58'
$ws.Range('E21').Value = '  -1.71%  '
$ws.Range('E22').Value = '  -1.57%  '
$ws.Range('E23').Value = '  -2.91%  '
$ws.Range('D24').Value = '1.999.28'
$ws.Range('E24').Value = '  -2.52%  '
$ws.Range('E25').Value = '  -3.80%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '151.54'
$ws.Range('E26').Value = '  -0.04%  '
$ws.Range('E27').Value = '  -1.67%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.153'
$ws.Range('E28').Value = '  -6.67%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.276'
$ws.Range('E29').Value = '  -1.98%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '115.79'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08933'
$ws.Range('E31').Value = '  +0.32%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.7588'
$ws.Range('E32').Value = '  -4.00%  '
$ws.Range('E33').Value = '  -3.84%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.454'
$ws.Range('E34').Value = '  -3.23%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.916'
$ws.Range('E35').Value = '  -0.05%  '
$ws.Range('E36').Value = '  -0.03%  '
$ws.Range('E37').Value = '  -0.47%  '
$ws.Range('E38').Value = '  -1.02%  '
$ws.Range('E39').Value = '  -0.85%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.929'
$ws.Range('E40').Value = '  +2.14%  '
$ws.Range('E41').Value = '  -1.83%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.373'
$ws.Range('E42').Value = '  +2.44%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.5288'
$ws.Range('E43').Value = '  -0.82%  '
$ws.Range('E44').Value = '  -2.92%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '8.507'
$ws.Range('E45').Value = '  -1.66%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.4998'
$ws.Range('E46').Value = '  -1.26%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.37'
$ws.Range('E47').Value = '  -2.61%  '
$ws.Range('B48').Value = 'PaxDollar'
$ws.Range('C48').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.000'
$ws.Range('E48').Value = '  -0.01%  '
$ws.Range('B49').Value = 'Quant'
$ws.Range('C49').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '103.89'
$ws.Range('E49').Value = '  -1.33%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.667'
$ws.Range('E50').Value = '  -1.87%  '
$ws.Range('E51').Value = '  -1.73%  '
